$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text for column B
$ws.Range("B1").Value = "Expiry Date (DD/MM/YYYY)"

# Remove the old placeholder row (###, DD/MM/YYYY) so the real ICCID/date
# row that used to be row 3 shifts up into row 2
$ws.Rows(2).Delete()

# Row 2 (previously row 3) now carries the ICCID / real date values already;
# give A2 the plain right-aligned numeric-style formatting used elsewhere in
# column A (picked up from an untouched cell further down the column) and set
# its final text value
$ws.Range("A1048576").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "8991102105546012952F "

# Update the expiry date value itself (keeps the dd/mm/yy formatting that
# shifted up from the old row 3)
$ws.Range("B2").Value = Get-Date -Year 2027 -Month 10 -Day 30 -Hour 0 -Minute 0 -Second 0

# Select B2 to match the saved view state
$ws.Range("B2").Select()
